# Update the EC (Estado de Cuenta) detail table on Hoja1.
#
# The original table (rows 16-34) listed every period for
# "IGNACIO TORRES BALSEIRO" (CC 73203467) followed by every period for
# "SANDY JAVIER DURANGO PEREGRINO" (CC 1143326442). The database was
# refreshed and Sandy's new records (periods 2205-2211) are now
# interleaved with Ignacio's rows by period, and Ignacio's own numbers
# were refreshed as well (new "Valor Mora" / F column amounts, plus an
# added 2112 period row at the top).
#
# Final target layout for rows 16-34 (columns C..G):
#   16: 73203467 | IGNACIO TORRES BALSEIRO       | 2112 | 10902 | 908526
#   17: 73203467 | IGNACIO TORRES BALSEIRO       | 2201 | 36341 | 908526
#   18: 73203467 | IGNACIO TORRES BALSEIRO       | 2202 | 36341 | 908526
#   19: 73203467 | IGNACIO TORRES BALSEIRO       | 2203 | 36341 | 908526
#   20: 73203467 | IGNACIO TORRES BALSEIRO       | 2204 | 36341 | 908526
#   21: 73203467 | IGNACIO TORRES BALSEIRO       | 2205 | 36341 | 908526
#   22: 1143326442 | SANDY JAVIER DURANGO PEREGRINO | 2205 | 47742 | 1193546
#   23: 73203467 | IGNACIO TORRES BALSEIRO       | 2206 | 36341 | 908526
#   24: 1143326442 | SANDY JAVIER DURANGO PEREGRINO | 2206 | 47742 | 1193546
#   25: 73203467 | IGNACIO TORRES BALSEIRO       | 2207 | 36341 | 908526
#   26: 1143326442 | SANDY JAVIER DURANGO PEREGRINO | 2207 | 47742 | 1193546
#   27: 73203467 | IGNACIO TORRES BALSEIRO       | 2208 | 36341 | 908526
#   28: 1143326442 | SANDY JAVIER DURANGO PEREGRINO | 2208 | 47742 | 1193546
#   29: 73203467 | IGNACIO TORRES BALSEIRO       | 2209 | 36341 | 908526
#   30: 1143326442 | SANDY JAVIER DURANGO PEREGRINO | 2209 | 47742 | 1193546
#   31: 73203467 | IGNACIO TORRES BALSEIRO       | 2210 | 36341 | 908526
#   32: 1143326442 | SANDY JAVIER DURANGO PEREGRINO | 2210 | 47742 | 1193546
#   33: 73203467 | IGNACIO TORRES BALSEIRO       | 2211 | 25439 | 908526
#   34: 1143326442 | SANDY JAVIER DURANGO PEREGRINO | 2211 | 33419 | 1193546

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ignacio = "IGNACIO TORRES BALSEIRO"
$ignacioId = "73203467"
$sandy = "SANDY JAVIER DURANGO PEREGRINO"
$sandyId = "1143326442"

$rows = @(
    @{ Row = 16; Id = $ignacioId; Name = $ignacio; Periodo = "2112"; Mora = 10902;  Salario = 908526 },
    @{ Row = 17; Id = $ignacioId; Name = $ignacio; Periodo = "2201"; Mora = 36341;  Salario = 908526 },
    @{ Row = 18; Id = $ignacioId; Name = $ignacio; Periodo = "2202"; Mora = 36341;  Salario = 908526 },
    @{ Row = 19; Id = $ignacioId; Name = $ignacio; Periodo = "2203"; Mora = 36341;  Salario = 908526 },
    @{ Row = 20; Id = $ignacioId; Name = $ignacio; Periodo = "2204"; Mora = 36341;  Salario = 908526 },
    @{ Row = 21; Id = $ignacioId; Name = $ignacio; Periodo = "2205"; Mora = 36341;  Salario = 908526 },
    @{ Row = 22; Id = $sandyId;   Name = $sandy;   Periodo = "2205"; Mora = 47742;  Salario = 1193546 },
    @{ Row = 23; Id = $ignacioId; Name = $ignacio; Periodo = "2206"; Mora = 36341;  Salario = 908526 },
    @{ Row = 24; Id = $sandyId;   Name = $sandy;   Periodo = "2206"; Mora = 47742;  Salario = 1193546 },
    @{ Row = 25; Id = $ignacioId; Name = $ignacio; Periodo = "2207"; Mora = 36341;  Salario = 908526 },
    @{ Row = 26; Id = $sandyId;   Name = $sandy;   Periodo = "2207"; Mora = 47742;  Salario = 1193546 },
    @{ Row = 27; Id = $ignacioId; Name = $ignacio; Periodo = "2208"; Mora = 36341;  Salario = 908526 },
    @{ Row = 28; Id = $sandyId;   Name = $sandy;   Periodo = "2208"; Mora = 47742;  Salario = 1193546 },
    @{ Row = 29; Id = $ignacioId; Name = $ignacio; Periodo = "2209"; Mora = 36341;  Salario = 908526 },
    @{ Row = 30; Id = $sandyId;   Name = $sandy;   Periodo = "2209"; Mora = 47742;  Salario = 1193546 },
    @{ Row = 31; Id = $ignacioId; Name = $ignacio; Periodo = "2210"; Mora = 36341;  Salario = 908526 },
    @{ Row = 32; Id = $sandyId;   Name = $sandy;   Periodo = "2210"; Mora = 47742;  Salario = 1193546 },
    @{ Row = 33; Id = $ignacioId; Name = $ignacio; Periodo = "2211"; Mora = 25439;  Salario = 908526 },
    @{ Row = 34; Id = $sandyId;   Name = $sandy;   Periodo = "2211"; Mora = 33419;  Salario = 1193546 }
)

foreach ($r in $rows) {
    $rowNum = $r.Row
    $ws.Range("C$rowNum").Value = $r.Id
    $ws.Range("D$rowNum").Value = $r.Name
    $ws.Range("E$rowNum").Value = $r.Periodo
    $ws.Range("F$rowNum").Value = $r.Mora
    $ws.Range("G$rowNum").Value = $r.Salario
}

Write-Host "EC table refreshed for rows 16-34"
